$d = $word.ActiveDocument

# Replace the date
$d.Content.Find.Execute("31-08-2023", $true, $false, $false, $false, $false, $true, 1, $false, "13-09-2023", 2)

# Replace all occurrences of the name "İrem" with "Mert"
$find = $d.Content.Find
$find.Text = "İrem"
$find.Replacement.Text = "Mert"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
